# 20230106 RENAME UI 코드
#
# The shared note in cell B3 of every sheet ("첫번째 단어는 대문자로
# 시작하고 (언더스코프)로 단어들 사이를 구분한다.") is replaced with the
# new naming-rule note ("카멜표기법 으로 구분한다."), i.e. the UI-code
# naming convention moved from underscore_case to camelCase.
#
# Sheet order in the workbook is: 환자정보 (1), 3DView (2), 3DView확대 (3).

$wb = $excel.ActiveWorkbook

$newNote = "카멜표기법 으로 구분한다."

$wsPatient = $wb.Worksheets.Item(1)   # 환자정보
$wsView    = $wb.Worksheets.Item(2)   # 3DView
$wsViewExp = $wb.Worksheets.Item(3)   # 3DView확대

# Update the merged B3:D3 note banner on every sheet with the new rule text.
$wsPatient.Range("B3").Value2 = $newNote
$wsView.Range("B3").Value2    = $newNote
$wsViewExp.Range("B3").Value2 = $newNote

# Restore the selections recorded for each sheet after the edit.
$wsPatient.Range("B3:D3").Select()
$wsView.Range("B3:D3").Select()
$wsViewExp.Range("F2").Select()
